# Pia Exercises unit 07 BDA Schol project KPIs
#
# Appends a new page of notes (22/02/2024 - "Codigos PAA" / "Cosdigos SQ")
# right after the paragraph "Objetivo: sacar tres tablas (como en el examen)".

$d = $word.ActiveDocument

# Locate the anchor paragraph (the last paragraph of the existing notes)
# and collapse the range to its end so we insert right after it, before
# the closing section properties.
$anchor = $d.Content
$found = $anchor.Find.Execute("Objetivo: sacar tres tablas (como en el examen)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph 'Objetivo: sacar tres tablas (como en el examen)' not found"
}

$insertionPoint = $d.Range($anchor.End, $anchor.End)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newContentXml = @"
<w:p $wNs><w:r><w:br w:type="page"/></w:r></w:p>
<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>22/02/2024</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>Codigos PAA</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>Los c&#243;digos PAA est&#225;n dividios en tres partes</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>PAA-1A1</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>1 -&gt; Linea</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>A -&gt; Objetivo</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>1 -&gt; Indicador</w:t></w:r></w:p>
<w:p $wNs/>
<w:p $wNs><w:r><w:t>Cosdigos SQ</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>Los c&#243;digos SQ</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>IPC</w:t></w:r><w:r><w:t>04.02</w:t></w:r></w:p>
<w:p $wNs><w:r><w:t>PC02 -&gt; Proceso</w:t></w:r></w:p>
"@

$insertionPoint.InsertXML($newContentXml)

Write-Host "Inserted 22/02/2024 PAA/SQ codes notes after the objetivo paragraph."
